$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append a new row (59) with the latest processed mail
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A59").Value = "Herinnering betaling"
$logs.Range("B59").Value = "mailmind.test@zohomail.eu"
$logs.Range("C59").Value = "Ik zie dat ik nog een openstaande betaling heb. Kunt u dit bevestigen?"
$logs.Range("D59").Value = "Factuur / Administratie"
$logs.Range("E59").Value = "Beste klant,`nBedankt voor uw bericht. Om u verder te kunnen helpen, heb ik wat meer informatie nodig om uw openstaande betaling te verifiëren. Kunt u mij alstublieft uw factuurnummer en/of klantnummer doorgeven? Zodra ik deze gegevens heb, zal ik direct voor u nakijken of er inderdaad nog een openstaande betaling is.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F59").Value = "2025-06-22 22:15:44"
$logs.Range("G59").Value = "Ja"

# Extend the conditional formatting ranges so they keep covering the
# "Categorie" (D) and "Beantwoord" (G) columns down to the new last row.
$dRules = $logs.Range("D2:D58").FormatConditions
$dRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D59"))

$gRules = $logs.Range("G2:G58").FormatConditions
$gRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G59"))

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet: the category counts/order shifted — Factuur /
#    Administratie moved up to #8 (unchanged count of 4 mails), pushing
#    Bestelling / Levering and Samenwerking / Partnerverzoek down, with
#    Samenwerking / Partnerverzoek's own count now at 4.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Factuur / Administratie"
$dash.Range("A9").Value = "Bestelling / Levering"
$dash.Range("A10").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B10").Value = 4
